# This script applies a row-permutation to the data rows (2-42) of the sheet.
# Columns D, L, M, N, O, P, R, S are shuffled between rows; row 14 is unchanged.
# We first snapshot every original value using .Value2 (resolves to the true scalar
# in this COM runtime, unlike .Value on some object types), then write the permuted
# values back using .Value so correct cell types / shared strings are produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orig = @{}
$orig[2] = @{ D = $ws.Range("D2").Value2; L = $ws.Range("L2").Value2; M = $ws.Range("M2").Value2; N = $ws.Range("N2").Value2; O = $ws.Range("O2").Value2; P = $ws.Range("P2").Value2; R = $ws.Range("R2").Value2; S = $ws.Range("S2").Value2 }
$orig[3] = @{ D = $ws.Range("D3").Value2; L = $ws.Range("L3").Value2; M = $ws.Range("M3").Value2; N = $ws.Range("N3").Value2; O = $ws.Range("O3").Value2; P = $ws.Range("P3").Value2; R = $ws.Range("R3").Value2; S = $ws.Range("S3").Value2 }
$orig[4] = @{ D = $ws.Range("D4").Value2; L = $ws.Range("L4").Value2; M = $ws.Range("M4").Value2; N = $ws.Range("N4").Value2; O = $ws.Range("O4").Value2; P = $ws.Range("P4").Value2; R = $ws.Range("R4").Value2; S = $ws.Range("S4").Value2 }
$orig[5] = @{ D = $ws.Range("D5").Value2; L = $ws.Range("L5").Value2; M = $ws.Range("M5").Value2; N = $ws.Range("N5").Value2; O = $ws.Range("O5").Value2; P = $ws.Range("P5").Value2; R = $ws.Range("R5").Value2; S = $ws.Range("S5").Value2 }
$orig[6] = @{ D = $ws.Range("D6").Value2; L = $ws.Range("L6").Value2; M = $ws.Range("M6").Value2; N = $ws.Range("N6").Value2; O = $ws.Range("O6").Value2; P = $ws.Range("P6").Value2; R = $ws.Range("R6").Value2; S = $ws.Range("S6").Value2 }
$orig[7] = @{ D = $ws.Range("D7").Value2; L = $ws.Range("L7").Value2; M = $ws.Range("M7").Value2; N = $ws.Range("N7").Value2; O = $ws.Range("O7").Value2; P = $ws.Range("P7").Value2; R = $ws.Range("R7").Value2; S = $ws.Range("S7").Value2 }
$orig[8] = @{ D = $ws.Range("D8").Value2; L = $ws.Range("L8").Value2; M = $ws.Range("M8").Value2; N = $ws.Range("N8").Value2; O = $ws.Range("O8").Value2; P = $ws.Range("P8").Value2; R = $ws.Range("R8").Value2; S = $ws.Range("S8").Value2 }
$orig[9] = @{ D = $ws.Range("D9").Value2; L = $ws.Range("L9").Value2; M = $ws.Range("M9").Value2; N = $ws.Range("N9").Value2; O = $ws.Range("O9").Value2; P = $ws.Range("P9").Value2; R = $ws.Range("R9").Value2; S = $ws.Range("S9").Value2 }
$orig[10] = @{ D = $ws.Range("D10").Value2; L = $ws.Range("L10").Value2; M = $ws.Range("M10").Value2; N = $ws.Range("N10").Value2; O = $ws.Range("O10").Value2; P = $ws.Range("P10").Value2; R = $ws.Range("R10").Value2; S = $ws.Range("S10").Value2 }
$orig[11] = @{ D = $ws.Range("D11").Value2; L = $ws.Range("L11").Value2; M = $ws.Range("M11").Value2; N = $ws.Range("N11").Value2; O = $ws.Range("O11").Value2; P = $ws.Range("P11").Value2; R = $ws.Range("R11").Value2; S = $ws.Range("S11").Value2 }
$orig[12] = @{ D = $ws.Range("D12").Value2; L = $ws.Range("L12").Value2; M = $ws.Range("M12").Value2; N = $ws.Range("N12").Value2; O = $ws.Range("O12").Value2; P = $ws.Range("P12").Value2; R = $ws.Range("R12").Value2; S = $ws.Range("S12").Value2 }
$orig[13] = @{ D = $ws.Range("D13").Value2; L = $ws.Range("L13").Value2; M = $ws.Range("M13").Value2; N = $ws.Range("N13").Value2; O = $ws.Range("O13").Value2; P = $ws.Range("P13").Value2; R = $ws.Range("R13").Value2; S = $ws.Range("S13").Value2 }
$orig[14] = @{ D = $ws.Range("D14").Value2; L = $ws.Range("L14").Value2; M = $ws.Range("M14").Value2; N = $ws.Range("N14").Value2; O = $ws.Range("O14").Value2; P = $ws.Range("P14").Value2; R = $ws.Range("R14").Value2; S = $ws.Range("S14").Value2 }
$orig[15] = @{ D = $ws.Range("D15").Value2; L = $ws.Range("L15").Value2; M = $ws.Range("M15").Value2; N = $ws.Range("N15").Value2; O = $ws.Range("O15").Value2; P = $ws.Range("P15").Value2; R = $ws.Range("R15").Value2; S = $ws.Range("S15").Value2 }
$orig[16] = @{ D = $ws.Range("D16").Value2; L = $ws.Range("L16").Value2; M = $ws.Range("M16").Value2; N = $ws.Range("N16").Value2; O = $ws.Range("O16").Value2; P = $ws.Range("P16").Value2; R = $ws.Range("R16").Value2; S = $ws.Range("S16").Value2 }
$orig[17] = @{ D = $ws.Range("D17").Value2; L = $ws.Range("L17").Value2; M = $ws.Range("M17").Value2; N = $ws.Range("N17").Value2; O = $ws.Range("O17").Value2; P = $ws.Range("P17").Value2; R = $ws.Range("R17").Value2; S = $ws.Range("S17").Value2 }
$orig[18] = @{ D = $ws.Range("D18").Value2; L = $ws.Range("L18").Value2; M = $ws.Range("M18").Value2; N = $ws.Range("N18").Value2; O = $ws.Range("O18").Value2; P = $ws.Range("P18").Value2; R = $ws.Range("R18").Value2; S = $ws.Range("S18").Value2 }
$orig[19] = @{ D = $ws.Range("D19").Value2; L = $ws.Range("L19").Value2; M = $ws.Range("M19").Value2; N = $ws.Range("N19").Value2; O = $ws.Range("O19").Value2; P = $ws.Range("P19").Value2; R = $ws.Range("R19").Value2; S = $ws.Range("S19").Value2 }
$orig[20] = @{ D = $ws.Range("D20").Value2; L = $ws.Range("L20").Value2; M = $ws.Range("M20").Value2; N = $ws.Range("N20").Value2; O = $ws.Range("O20").Value2; P = $ws.Range("P20").Value2; R = $ws.Range("R20").Value2; S = $ws.Range("S20").Value2 }
$orig[21] = @{ D = $ws.Range("D21").Value2; L = $ws.Range("L21").Value2; M = $ws.Range("M21").Value2; N = $ws.Range("N21").Value2; O = $ws.Range("O21").Value2; P = $ws.Range("P21").Value2; R = $ws.Range("R21").Value2; S = $ws.Range("S21").Value2 }
$orig[22] = @{ D = $ws.Range("D22").Value2; L = $ws.Range("L22").Value2; M = $ws.Range("M22").Value2; N = $ws.Range("N22").Value2; O = $ws.Range("O22").Value2; P = $ws.Range("P22").Value2; R = $ws.Range("R22").Value2; S = $ws.Range("S22").Value2 }
$orig[23] = @{ D = $ws.Range("D23").Value2; L = $ws.Range("L23").Value2; M = $ws.Range("M23").Value2; N = $ws.Range("N23").Value2; O = $ws.Range("O23").Value2; P = $ws.Range("P23").Value2; R = $ws.Range("R23").Value2; S = $ws.Range("S23").Value2 }
$orig[24] = @{ D = $ws.Range("D24").Value2; L = $ws.Range("L24").Value2; M = $ws.Range("M24").Value2; N = $ws.Range("N24").Value2; O = $ws.Range("O24").Value2; P = $ws.Range("P24").Value2; R = $ws.Range("R24").Value2; S = $ws.Range("S24").Value2 }
$orig[25] = @{ D = $ws.Range("D25").Value2; L = $ws.Range("L25").Value2; M = $ws.Range("M25").Value2; N = $ws.Range("N25").Value2; O = $ws.Range("O25").Value2; P = $ws.Range("P25").Value2; R = $ws.Range("R25").Value2; S = $ws.Range("S25").Value2 }
$orig[26] = @{ D = $ws.Range("D26").Value2; L = $ws.Range("L26").Value2; M = $ws.Range("M26").Value2; N = $ws.Range("N26").Value2; O = $ws.Range("O26").Value2; P = $ws.Range("P26").Value2; R = $ws.Range("R26").Value2; S = $ws.Range("S26").Value2 }
$orig[27] = @{ D = $ws.Range("D27").Value2; L = $ws.Range("L27").Value2; M = $ws.Range("M27").Value2; N = $ws.Range("N27").Value2; O = $ws.Range("O27").Value2; P = $ws.Range("P27").Value2; R = $ws.Range("R27").Value2; S = $ws.Range("S27").Value2 }
$orig[28] = @{ D = $ws.Range("D28").Value2; L = $ws.Range("L28").Value2; M = $ws.Range("M28").Value2; N = $ws.Range("N28").Value2; O = $ws.Range("O28").Value2; P = $ws.Range("P28").Value2; R = $ws.Range("R28").Value2; S = $ws.Range("S28").Value2 }
$orig[29] = @{ D = $ws.Range("D29").Value2; L = $ws.Range("L29").Value2; M = $ws.Range("M29").Value2; N = $ws.Range("N29").Value2; O = $ws.Range("O29").Value2; P = $ws.Range("P29").Value2; R = $ws.Range("R29").Value2; S = $ws.Range("S29").Value2 }
$orig[30] = @{ D = $ws.Range("D30").Value2; L = $ws.Range("L30").Value2; M = $ws.Range("M30").Value2; N = $ws.Range("N30").Value2; O = $ws.Range("O30").Value2; P = $ws.Range("P30").Value2; R = $ws.Range("R30").Value2; S = $ws.Range("S30").Value2 }
$orig[31] = @{ D = $ws.Range("D31").Value2; L = $ws.Range("L31").Value2; M = $ws.Range("M31").Value2; N = $ws.Range("N31").Value2; O = $ws.Range("O31").Value2; P = $ws.Range("P31").Value2; R = $ws.Range("R31").Value2; S = $ws.Range("S31").Value2 }
$orig[32] = @{ D = $ws.Range("D32").Value2; L = $ws.Range("L32").Value2; M = $ws.Range("M32").Value2; N = $ws.Range("N32").Value2; O = $ws.Range("O32").Value2; P = $ws.Range("P32").Value2; R = $ws.Range("R32").Value2; S = $ws.Range("S32").Value2 }
$orig[33] = @{ D = $ws.Range("D33").Value2; L = $ws.Range("L33").Value2; M = $ws.Range("M33").Value2; N = $ws.Range("N33").Value2; O = $ws.Range("O33").Value2; P = $ws.Range("P33").Value2; R = $ws.Range("R33").Value2; S = $ws.Range("S33").Value2 }
$orig[34] = @{ D = $ws.Range("D34").Value2; L = $ws.Range("L34").Value2; M = $ws.Range("M34").Value2; N = $ws.Range("N34").Value2; O = $ws.Range("O34").Value2; P = $ws.Range("P34").Value2; R = $ws.Range("R34").Value2; S = $ws.Range("S34").Value2 }
$orig[35] = @{ D = $ws.Range("D35").Value2; L = $ws.Range("L35").Value2; M = $ws.Range("M35").Value2; N = $ws.Range("N35").Value2; O = $ws.Range("O35").Value2; P = $ws.Range("P35").Value2; R = $ws.Range("R35").Value2; S = $ws.Range("S35").Value2 }
$orig[36] = @{ D = $ws.Range("D36").Value2; L = $ws.Range("L36").Value2; M = $ws.Range("M36").Value2; N = $ws.Range("N36").Value2; O = $ws.Range("O36").Value2; P = $ws.Range("P36").Value2; R = $ws.Range("R36").Value2; S = $ws.Range("S36").Value2 }
$orig[37] = @{ D = $ws.Range("D37").Value2; L = $ws.Range("L37").Value2; M = $ws.Range("M37").Value2; N = $ws.Range("N37").Value2; O = $ws.Range("O37").Value2; P = $ws.Range("P37").Value2; R = $ws.Range("R37").Value2; S = $ws.Range("S37").Value2 }
$orig[38] = @{ D = $ws.Range("D38").Value2; L = $ws.Range("L38").Value2; M = $ws.Range("M38").Value2; N = $ws.Range("N38").Value2; O = $ws.Range("O38").Value2; P = $ws.Range("P38").Value2; R = $ws.Range("R38").Value2; S = $ws.Range("S38").Value2 }
$orig[39] = @{ D = $ws.Range("D39").Value2; L = $ws.Range("L39").Value2; M = $ws.Range("M39").Value2; N = $ws.Range("N39").Value2; O = $ws.Range("O39").Value2; P = $ws.Range("P39").Value2; R = $ws.Range("R39").Value2; S = $ws.Range("S39").Value2 }
$orig[40] = @{ D = $ws.Range("D40").Value2; L = $ws.Range("L40").Value2; M = $ws.Range("M40").Value2; N = $ws.Range("N40").Value2; O = $ws.Range("O40").Value2; P = $ws.Range("P40").Value2; R = $ws.Range("R40").Value2; S = $ws.Range("S40").Value2 }
$orig[41] = @{ D = $ws.Range("D41").Value2; L = $ws.Range("L41").Value2; M = $ws.Range("M41").Value2; N = $ws.Range("N41").Value2; O = $ws.Range("O41").Value2; P = $ws.Range("P41").Value2; R = $ws.Range("R41").Value2; S = $ws.Range("S41").Value2 }
$orig[42] = @{ D = $ws.Range("D42").Value2; L = $ws.Range("L42").Value2; M = $ws.Range("M42").Value2; N = $ws.Range("N42").Value2; O = $ws.Range("O42").Value2; P = $ws.Range("P42").Value2; R = $ws.Range("R42").Value2; S = $ws.Range("S42").Value2 }

# Write back the permuted values (row $dst now holds what used to be in row $src)
$ws.Range("D2").Value = $orig[36].D
$ws.Range("L2").Value = $orig[36].L
$ws.Range("M2").Value = $orig[36].M
$ws.Range("N2").Value = $orig[36].N
$ws.Range("O2").Value = $orig[36].O
$ws.Range("P2").Value = $orig[36].P
$ws.Range("R2").Value = $orig[36].R
$ws.Range("S2").Value = $orig[36].S
$ws.Range("D3").Value = $orig[37].D
$ws.Range("L3").Value = $orig[37].L
$ws.Range("M3").Value = $orig[37].M
$ws.Range("N3").Value = $orig[37].N
$ws.Range("O3").Value = $orig[37].O
$ws.Range("P3").Value = $orig[37].P
$ws.Range("R3").Value = $orig[37].R
$ws.Range("S3").Value = $orig[37].S
$ws.Range("D4").Value = $orig[24].D
$ws.Range("L4").Value = $orig[24].L
$ws.Range("M4").Value = $orig[24].M
$ws.Range("N4").Value = $orig[24].N
$ws.Range("O4").Value = $orig[24].O
$ws.Range("P4").Value = $orig[24].P
$ws.Range("R4").Value = $orig[24].R
$ws.Range("S4").Value = $orig[24].S
$ws.Range("D5").Value = $orig[31].D
$ws.Range("L5").Value = $orig[31].L
$ws.Range("M5").Value = $orig[31].M
$ws.Range("N5").Value = $orig[31].N
$ws.Range("O5").Value = $orig[31].O
$ws.Range("P5").Value = $orig[31].P
$ws.Range("R5").Value = $orig[31].R
$ws.Range("S5").Value = $orig[31].S
$ws.Range("D6").Value = $orig[32].D
$ws.Range("L6").Value = $orig[32].L
$ws.Range("M6").Value = $orig[32].M
$ws.Range("N6").Value = $orig[32].N
$ws.Range("O6").Value = $orig[32].O
$ws.Range("P6").Value = $orig[32].P
$ws.Range("R6").Value = $orig[32].R
$ws.Range("S6").Value = $orig[32].S
$ws.Range("D7").Value = $orig[25].D
$ws.Range("L7").Value = $orig[25].L
$ws.Range("M7").Value = $orig[25].M
$ws.Range("N7").Value = $orig[25].N
$ws.Range("O7").Value = $orig[25].O
$ws.Range("P7").Value = $orig[25].P
$ws.Range("R7").Value = $orig[25].R
$ws.Range("S7").Value = $orig[25].S
$ws.Range("D8").Value = $orig[26].D
$ws.Range("L8").Value = $orig[26].L
$ws.Range("M8").Value = $orig[26].M
$ws.Range("N8").Value = $orig[26].N
$ws.Range("O8").Value = $orig[26].O
$ws.Range("P8").Value = $orig[26].P
$ws.Range("R8").Value = $orig[26].R
$ws.Range("S8").Value = $orig[26].S
$ws.Range("D9").Value = $orig[42].D
$ws.Range("L9").Value = $orig[42].L
$ws.Range("M9").Value = $orig[42].M
$ws.Range("N9").Value = $orig[42].N
$ws.Range("O9").Value = $orig[42].O
$ws.Range("P9").Value = $orig[42].P
$ws.Range("R9").Value = $orig[42].R
$ws.Range("S9").Value = $orig[42].S
$ws.Range("D10").Value = $orig[16].D
$ws.Range("L10").Value = $orig[16].L
$ws.Range("M10").Value = $orig[16].M
$ws.Range("N10").Value = $orig[16].N
$ws.Range("O10").Value = $orig[16].O
$ws.Range("P10").Value = $orig[16].P
$ws.Range("R10").Value = $orig[16].R
$ws.Range("S10").Value = $orig[16].S
$ws.Range("D11").Value = $orig[17].D
$ws.Range("L11").Value = $orig[17].L
$ws.Range("M11").Value = $orig[17].M
$ws.Range("N11").Value = $orig[17].N
$ws.Range("O11").Value = $orig[17].O
$ws.Range("P11").Value = $orig[17].P
$ws.Range("R11").Value = $orig[17].R
$ws.Range("S11").Value = $orig[17].S
$ws.Range("D12").Value = $orig[21].D
$ws.Range("L12").Value = $orig[21].L
$ws.Range("M12").Value = $orig[21].M
$ws.Range("N12").Value = $orig[21].N
$ws.Range("O12").Value = $orig[21].O
$ws.Range("P12").Value = $orig[21].P
$ws.Range("R12").Value = $orig[21].R
$ws.Range("S12").Value = $orig[21].S
$ws.Range("D13").Value = $orig[40].D
$ws.Range("L13").Value = $orig[40].L
$ws.Range("M13").Value = $orig[40].M
$ws.Range("N13").Value = $orig[40].N
$ws.Range("O13").Value = $orig[40].O
$ws.Range("P13").Value = $orig[40].P
$ws.Range("R13").Value = $orig[40].R
$ws.Range("S13").Value = $orig[40].S
$ws.Range("D15").Value = $orig[33].D
$ws.Range("L15").Value = $orig[33].L
$ws.Range("M15").Value = $orig[33].M
$ws.Range("N15").Value = $orig[33].N
$ws.Range("O15").Value = $orig[33].O
$ws.Range("P15").Value = $orig[33].P
$ws.Range("R15").Value = $orig[33].R
$ws.Range("S15").Value = $orig[33].S
$ws.Range("D16").Value = $orig[22].D
$ws.Range("L16").Value = $orig[22].L
$ws.Range("M16").Value = $orig[22].M
$ws.Range("N16").Value = $orig[22].N
$ws.Range("O16").Value = $orig[22].O
$ws.Range("P16").Value = $orig[22].P
$ws.Range("R16").Value = $orig[22].R
$ws.Range("S16").Value = $orig[22].S
$ws.Range("D17").Value = $orig[41].D
$ws.Range("L17").Value = $orig[41].L
$ws.Range("M17").Value = $orig[41].M
$ws.Range("N17").Value = $orig[41].N
$ws.Range("O17").Value = $orig[41].O
$ws.Range("P17").Value = $orig[41].P
$ws.Range("R17").Value = $orig[41].R
$ws.Range("S17").Value = $orig[41].S
$ws.Range("D18").Value = $orig[38].D
$ws.Range("L18").Value = $orig[38].L
$ws.Range("M18").Value = $orig[38].M
$ws.Range("N18").Value = $orig[38].N
$ws.Range("O18").Value = $orig[38].O
$ws.Range("P18").Value = $orig[38].P
$ws.Range("R18").Value = $orig[38].R
$ws.Range("S18").Value = $orig[38].S
$ws.Range("D19").Value = $orig[39].D
$ws.Range("L19").Value = $orig[39].L
$ws.Range("M19").Value = $orig[39].M
$ws.Range("N19").Value = $orig[39].N
$ws.Range("O19").Value = $orig[39].O
$ws.Range("P19").Value = $orig[39].P
$ws.Range("R19").Value = $orig[39].R
$ws.Range("S19").Value = $orig[39].S
$ws.Range("D20").Value = $orig[4].D
$ws.Range("L20").Value = $orig[4].L
$ws.Range("M20").Value = $orig[4].M
$ws.Range("N20").Value = $orig[4].N
$ws.Range("O20").Value = $orig[4].O
$ws.Range("P20").Value = $orig[4].P
$ws.Range("R20").Value = $orig[4].R
$ws.Range("S20").Value = $orig[4].S
$ws.Range("D21").Value = $orig[11].D
$ws.Range("L21").Value = $orig[11].L
$ws.Range("M21").Value = $orig[11].M
$ws.Range("N21").Value = $orig[11].N
$ws.Range("O21").Value = $orig[11].O
$ws.Range("P21").Value = $orig[11].P
$ws.Range("R21").Value = $orig[11].R
$ws.Range("S21").Value = $orig[11].S
$ws.Range("D22").Value = $orig[12].D
$ws.Range("L22").Value = $orig[12].L
$ws.Range("M22").Value = $orig[12].M
$ws.Range("N22").Value = $orig[12].N
$ws.Range("O22").Value = $orig[12].O
$ws.Range("P22").Value = $orig[12].P
$ws.Range("R22").Value = $orig[12].R
$ws.Range("S22").Value = $orig[12].S
$ws.Range("D23").Value = $orig[9].D
$ws.Range("L23").Value = $orig[9].L
$ws.Range("M23").Value = $orig[9].M
$ws.Range("N23").Value = $orig[9].N
$ws.Range("O23").Value = $orig[9].O
$ws.Range("P23").Value = $orig[9].P
$ws.Range("R23").Value = $orig[9].R
$ws.Range("S23").Value = $orig[9].S
$ws.Range("D24").Value = $orig[10].D
$ws.Range("L24").Value = $orig[10].L
$ws.Range("M24").Value = $orig[10].M
$ws.Range("N24").Value = $orig[10].N
$ws.Range("O24").Value = $orig[10].O
$ws.Range("P24").Value = $orig[10].P
$ws.Range("R24").Value = $orig[10].R
$ws.Range("S24").Value = $orig[10].S
$ws.Range("D25").Value = $orig[28].D
$ws.Range("L25").Value = $orig[28].L
$ws.Range("M25").Value = $orig[28].M
$ws.Range("N25").Value = $orig[28].N
$ws.Range("O25").Value = $orig[28].O
$ws.Range("P25").Value = $orig[28].P
$ws.Range("R25").Value = $orig[28].R
$ws.Range("S25").Value = $orig[28].S
$ws.Range("D26").Value = $orig[29].D
$ws.Range("L26").Value = $orig[29].L
$ws.Range("M26").Value = $orig[29].M
$ws.Range("N26").Value = $orig[29].N
$ws.Range("O26").Value = $orig[29].O
$ws.Range("P26").Value = $orig[29].P
$ws.Range("R26").Value = $orig[29].R
$ws.Range("S26").Value = $orig[29].S
$ws.Range("D27").Value = $orig[18].D
$ws.Range("L27").Value = $orig[18].L
$ws.Range("M27").Value = $orig[18].M
$ws.Range("N27").Value = $orig[18].N
$ws.Range("O27").Value = $orig[18].O
$ws.Range("P27").Value = $orig[18].P
$ws.Range("R27").Value = $orig[18].R
$ws.Range("S27").Value = $orig[18].S
$ws.Range("D28").Value = $orig[19].D
$ws.Range("L28").Value = $orig[19].L
$ws.Range("M28").Value = $orig[19].M
$ws.Range("N28").Value = $orig[19].N
$ws.Range("O28").Value = $orig[19].O
$ws.Range("P28").Value = $orig[19].P
$ws.Range("R28").Value = $orig[19].R
$ws.Range("S28").Value = $orig[19].S
$ws.Range("D29").Value = $orig[27].D
$ws.Range("L29").Value = $orig[27].L
$ws.Range("M29").Value = $orig[27].M
$ws.Range("N29").Value = $orig[27].N
$ws.Range("O29").Value = $orig[27].O
$ws.Range("P29").Value = $orig[27].P
$ws.Range("R29").Value = $orig[27].R
$ws.Range("S29").Value = $orig[27].S
$ws.Range("D30").Value = $orig[8].D
$ws.Range("L30").Value = $orig[8].L
$ws.Range("M30").Value = $orig[8].M
$ws.Range("N30").Value = $orig[8].N
$ws.Range("O30").Value = $orig[8].O
$ws.Range("P30").Value = $orig[8].P
$ws.Range("R30").Value = $orig[8].R
$ws.Range("S30").Value = $orig[8].S
$ws.Range("D31").Value = $orig[34].D
$ws.Range("L31").Value = $orig[34].L
$ws.Range("M31").Value = $orig[34].M
$ws.Range("N31").Value = $orig[34].N
$ws.Range("O31").Value = $orig[34].O
$ws.Range("P31").Value = $orig[34].P
$ws.Range("R31").Value = $orig[34].R
$ws.Range("S31").Value = $orig[34].S
$ws.Range("D32").Value = $orig[35].D
$ws.Range("L32").Value = $orig[35].L
$ws.Range("M32").Value = $orig[35].M
$ws.Range("N32").Value = $orig[35].N
$ws.Range("O32").Value = $orig[35].O
$ws.Range("P32").Value = $orig[35].P
$ws.Range("R32").Value = $orig[35].R
$ws.Range("S32").Value = $orig[35].S
$ws.Range("D33").Value = $orig[13].D
$ws.Range("L33").Value = $orig[13].L
$ws.Range("M33").Value = $orig[13].M
$ws.Range("N33").Value = $orig[13].N
$ws.Range("O33").Value = $orig[13].O
$ws.Range("P33").Value = $orig[13].P
$ws.Range("R33").Value = $orig[13].R
$ws.Range("S33").Value = $orig[13].S
$ws.Range("D34").Value = $orig[23].D
$ws.Range("L34").Value = $orig[23].L
$ws.Range("M34").Value = $orig[23].M
$ws.Range("N34").Value = $orig[23].N
$ws.Range("O34").Value = $orig[23].O
$ws.Range("P34").Value = $orig[23].P
$ws.Range("R34").Value = $orig[23].R
$ws.Range("S34").Value = $orig[23].S
$ws.Range("D35").Value = $orig[5].D
$ws.Range("L35").Value = $orig[5].L
$ws.Range("M35").Value = $orig[5].M
$ws.Range("N35").Value = $orig[5].N
$ws.Range("O35").Value = $orig[5].O
$ws.Range("P35").Value = $orig[5].P
$ws.Range("R35").Value = $orig[5].R
$ws.Range("S35").Value = $orig[5].S
$ws.Range("D36").Value = $orig[6].D
$ws.Range("L36").Value = $orig[6].L
$ws.Range("M36").Value = $orig[6].M
$ws.Range("N36").Value = $orig[6].N
$ws.Range("O36").Value = $orig[6].O
$ws.Range("P36").Value = $orig[6].P
$ws.Range("R36").Value = $orig[6].R
$ws.Range("S36").Value = $orig[6].S
$ws.Range("D37").Value = $orig[7].D
$ws.Range("L37").Value = $orig[7].L
$ws.Range("M37").Value = $orig[7].M
$ws.Range("N37").Value = $orig[7].N
$ws.Range("O37").Value = $orig[7].O
$ws.Range("P37").Value = $orig[7].P
$ws.Range("R37").Value = $orig[7].R
$ws.Range("S37").Value = $orig[7].S
$ws.Range("D38").Value = $orig[30].D
$ws.Range("L38").Value = $orig[30].L
$ws.Range("M38").Value = $orig[30].M
$ws.Range("N38").Value = $orig[30].N
$ws.Range("O38").Value = $orig[30].O
$ws.Range("P38").Value = $orig[30].P
$ws.Range("R38").Value = $orig[30].R
$ws.Range("S38").Value = $orig[30].S
$ws.Range("D39").Value = $orig[15].D
$ws.Range("L39").Value = $orig[15].L
$ws.Range("M39").Value = $orig[15].M
$ws.Range("N39").Value = $orig[15].N
$ws.Range("O39").Value = $orig[15].O
$ws.Range("P39").Value = $orig[15].P
$ws.Range("R39").Value = $orig[15].R
$ws.Range("S39").Value = $orig[15].S
$ws.Range("D40").Value = $orig[2].D
$ws.Range("L40").Value = $orig[2].L
$ws.Range("M40").Value = $orig[2].M
$ws.Range("N40").Value = $orig[2].N
$ws.Range("O40").Value = $orig[2].O
$ws.Range("P40").Value = $orig[2].P
$ws.Range("R40").Value = $orig[2].R
$ws.Range("S40").Value = $orig[2].S
$ws.Range("D41").Value = $orig[3].D
$ws.Range("L41").Value = $orig[3].L
$ws.Range("M41").Value = $orig[3].M
$ws.Range("N41").Value = $orig[3].N
$ws.Range("O41").Value = $orig[3].O
$ws.Range("P41").Value = $orig[3].P
$ws.Range("R41").Value = $orig[3].R
$ws.Range("S41").Value = $orig[3].S
$ws.Range("D42").Value = $orig[20].D
$ws.Range("L42").Value = $orig[20].L
$ws.Range("M42").Value = $orig[20].M
$ws.Range("N42").Value = $orig[20].N
$ws.Range("O42").Value = $orig[20].O
$ws.Range("P42").Value = $orig[20].P
$ws.Range("R42").Value = $orig[20].R
$ws.Range("S42").Value = $orig[20].S
